$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 54
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0.2
$ws.Range("D2").Value = 0.2
$ws.Range("E2").Value = 0.2
$ws.Range("F2").Value = 1.1
